# GMS Data Release 1
# Rename the "patient_id" row entry to "participant_id" and move the
# active selection to F9 (matching the authored diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B8 held "patient_id" -> rename to "participant_id"
$ws.Range("B8").Value = "participant_id"

# Move the active cell/selection to F9
[void]$ws.Range("F9").Select()
